$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for the new days (20 and 21) of May/2025 right before
# the start of the April/2025 block (old row 21), pushing everything else
# down by two rows.
$ws.Rows("21:22").Insert()

# Day 20 - Maio/2025
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = 33090.3
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 2025
$ws.Range("E21").Value = "05/2025"

# Day 21 - Maio/2025
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = 24686.66
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 2025
$ws.Range("E22").Value = "05/2025"
